# Natmi following Dr Hou advice
# Rebuild the Wnt1-Fzd9 LR-pair sheet: the original data only covered the
# "FAPs" sending cluster; the revised numbers add an "ECs" sending-cluster
# block (new rows 2-6) and recompute the existing "FAPs" block (now rows
# 7-11) against the full 3-sample background.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs","Wnt1","Fzd9","ECs",1,0.3333333333333333,0.01948966666666667,0.058469,0.07096062449330311,0.07096062449330311,3,1,0.3962276666666666,1.188683,0.1703050239984229,0.1960950025990754,0.007722345147444443,0.06950110632699999,0.01208495085727506,0.01391502384444629),
    @("ECs","Wnt1","Fzd9","FAPs",1,0.3333333333333333,0.01948966666666667,0.058469,0.07096062449330311,0.07096062449330311,3,1,0.6295006666666667,1.888502,0.270569511325618,0.3115429467724861,0.01226875815977778,0.110418823438,0.0191997814925137,0.02210728205945951),
    @("ECs","Wnt1","Fzd9","M1",1,0.3333333333333333,0.01948966666666667,0.058469,0.07096062449330311,0.07096062449330311,3,1,0.1632393333333333,0.489718,0.07016289098309612,0.08078794134585419,0.003181480193555555,0.028633321742,0.004978802560416046,0.005732762769430156),
    @("ECs","Wnt1","Fzd9","M2",1,0.3333333333333333,0.01948966666666667,0.058469,0.07096062449330311,0.07096062449330311,3,1,0.2196503333333334,0.6589510000000001,0.0944092460889781,0.1087060200723518,0.00428091177988889,0.03852820601900001,0.00669933906041582,0.007713847070515627),
    @("ECs","Wnt1","Fzd9","sCs",1,0.3333333333333333,0.01948966666666667,0.058469,0.07096062449330311,0.07096062449330311,2,1,0.9179585,1.835917,0.394553327603885,0.3028680892102324,0.01789070517883333,0.107344231073,0.02799775052268248,0.02149170874945153),
    @("FAPs","Wnt1","Fzd9","ECs",3,1,0.255165,0.7654949999999999,0.9290393755066968,0.9290393755066969,3,1,0.3962276666666666,1.188683,0.1703050239984229,0.1960950025990754,0.101103432565,0.9099308930849999,0.1582200731411478,0.1821799787546291),
    @("FAPs","Wnt1","Fzd9","FAPs",3,1,0.255165,0.7654949999999999,0.9290393755066968,0.9290393755066969,3,1,0.6295006666666667,1.888502,0.270569511325618,0.3115429467724861,0.16062653761,1.44563883849,0.2513697298331043,0.2894356647130266),
    @("FAPs","Wnt1","Fzd9","M1",3,1,0.255165,0.7654949999999999,0.9290393755066968,0.9290393755066969,3,1,0.1632393333333333,0.489718,0.07016289098309612,0.08078794134585419,0.04165296448999999,0.3748766804099999,0.06518408842268007,0.07505517857642403),
    @("FAPs","Wnt1","Fzd9","M2",3,1,0.255165,0.7654949999999999,0.9290393755066968,0.9290393755066969,3,1,0.2196503333333334,0.6589510000000001,0.0944092460889781,0.1087060200723518,0.05604707730500001,0.504423695745,0.08770990702856228,0.1009921730018362),
    @("FAPs","Wnt1","Fzd9","sCs",3,1,0.255165,0.7654949999999999,0.9290393755066968,0.9290393755066969,2,1,0.9179585,1.835917,0.394553327603885,0.3028680892102324,0.2342308806525,1.405385283915,0.3665555770812025,0.2813763804607809)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
